$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.20429999999999
$ws.Range("A9").Value = -21.44649999999998
$ws.Range("C9").Value = -10.2711
$ws.Range("C11").Value = -13.13960000000001
$ws.Range("A18").Value = -22.23110000000002
$ws.Range("A20").Value = -21.42369999999999
$ws.Range("C23").Value = -12.0717
$ws.Range("C24").Value = -12.7854
$ws.Range("C26").Value = -12.64200000000001
$ws.Range("A27").Value = -21.90829999999999
$ws.Range("C34").Value = -12.19720000000001
$ws.Range("A35").Value = -21.56679999999998
$ws.Range("C35").Value = -12.57440000000001
$ws.Range("C48").Value = -11.68639999999999
$ws.Range("C49").Value = -13.59109999999999
$ws.Range("C52").Value = -10.7517
$ws.Range("C66").Value = -11.1536
$ws.Range("C67").Value = -11.2998
$ws.Range("A69").Value = -21.49189999999997
$ws.Range("A76").Value = -19.38139999999999
$ws.Range("A78").Value = -19.83159999999998
$ws.Range("C78").Value = -13.2926
$ws.Range("C80").Value = -13.11340000000001
$ws.Range("A82").Value = -21.787
$ws.Range("A83").Value = -21.56519999999999
$ws.Range("A93").Value = -21.42600000000001
$ws.Range("C99").Value = -12.8319
$ws.Range("C104").Value = -12.90790000000001
